# Generate Report for handoff
#
# A new handoff package was produced for "b.md.md". This updates the
# status/handoff information for that row on the Overview, zh-cn and
# de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
# Row 3 corresponds to b.md.md (A3). Columns B (zh-cn) and C (de-de)
# hold the localization status for that file.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# --- zh-cn sheet ------------------------------------------------------
# Row 3 corresponds to b.md.md (A3). Update Status, Latest Handoff File
# and Latest Handoff Datetime, and fix up the hyperlink display text so
# it matches the new handoff file name.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-19 04:20:32"
foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# --- de-de sheet ------------------------------------------------------
# Row 3 corresponds to b.md.md (A3). Update Status, Latest Handoff File
# and Latest Handoff Datetime, and fix up the hyperlink display text so
# it matches the new handoff file name.
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-19 04:20:41"
foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
